# The commit removes the single reviewer comment that was attached to the
# end of the "Zakljucak" (Conclusion) paragraph discussing the LegalRuleML
# parser contribution. Deleting the Comment object removes both the
# <w:commentReference> run in word/document.xml and the comment entry
# itself from word/comments.xml.

$d = $word.ActiveDocument

if ($d.Comments.Count -gt 0) {
    # Walk backwards in case deleting shifts indices.
    for ($i = $d.Comments.Count; $i -ge 1; $i--) {
        $d.Comments.Item($i).Delete()
    }
}
